# Auto-generated script to apply scheduled market-data refresh to Leve profit tables
# Updates currentAveragePrice / NQ / HQ price & profit columns (H:N) for the affected leve rows
# across the ALC, ARM, BSM, CRP, CUL, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 4
$ws.Range("H4").Value = 2071.4
$ws.Range("I4").Value = 2071.4
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2071.4
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -1957.4
$ws.Range("N4").ClearContents()

# Row 33
$ws.Range("H33").Value = 325672.22
$ws.Range("I33").Value = 518.2857
$ws.Range("K33").Value = 518.2857
$ws.Range("M33").Value = -289.2857

# Row 39
$ws.Range("H39").Value = 1133947.1
$ws.Range("I39").Value = 1443078.1
$ws.Range("J39").Value = 466.66666
$ws.Range("K39").Value = 4329234.300000001
$ws.Range("L39").Value = 1399.99998
$ws.Range("M39").Value = -4328938.300000001
$ws.Range("N39").Value = -1991.99998

# Row 44
$ws.Range("H44").Value = 17000
$ws.Range("J44").Value = 17000
$ws.Range("L44").Value = 17000
$ws.Range("N44").Value = -17924

# Row 69
$ws.Range("H69").Value = 3515
$ws.Range("I69").Value = 4266.6665
$ws.Range("J69").Value = 3264.4443
$ws.Range("K69").Value = 12799.9995
$ws.Range("L69").Value = 9793.332900000001
$ws.Range("M69").Value = -11925.9995
$ws.Range("N69").Value = -11541.3329

# Row 72
$ws.Range("H72").Value = 3515
$ws.Range("I72").Value = 4266.6665
$ws.Range("J72").Value = 3264.4443
$ws.Range("K72").Value = 38399.9985
$ws.Range("L72").Value = 29379.9987
$ws.Range("M72").Value = -34031.9985
$ws.Range("N72").Value = -38115.9987

# Row 129
$ws.Range("H129").Value = 3482.3684
$ws.Range("J129").Value = 1104.9048
$ws.Range("L129").Value = 3314.7144
$ws.Range("N129").Value = -13314.7144

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 36713.785
$ws.Range("I2").Value = 1090.72
$ws.Range("J2").Value = 333572.66
$ws.Range("K2").Value = 1090.72
$ws.Range("L2").Value = 333572.66
$ws.Range("M2").Value = -977.72
$ws.Range("N2").Value = -333798.66

# Row 32
$ws.Range("H32").Value = 5187.66
$ws.Range("I32").Value = 4561.3403
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 4561.3403
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -4274.3403
$ws.Range("N32").Value = -15574

# Row 61
$ws.Range("H61").Value = 2786.5
$ws.Range("I61").Value = 1750
$ws.Range("K61").Value = 1750
$ws.Range("M61").Value = -1538

# Row 74
$ws.Range("H74").Value = 873.4054
$ws.Range("I74").Value = 808.1
$ws.Range("K74").Value = 808.1
$ws.Range("M74").Value = 65.89999999999998

# Row 77
$ws.Range("H77").Value = 873.4054
$ws.Range("I77").Value = 808.1
$ws.Range("K77").Value = 4040.5
$ws.Range("M77").Value = 327.5

# Row 102
$ws.Range("H102").Value = 73952.71000000001
$ws.Range("I102").Value = 127734.75
$ws.Range("J102").Value = 2243.3333
$ws.Range("K102").Value = 127734.75
$ws.Range("L102").Value = 2243.3333
$ws.Range("M102").Value = -126112.75
$ws.Range("N102").Value = -5487.3333

# Row 111
$ws.Range("H111").Value = 32999.668
$ws.Range("J111").Value = 32999.668
$ws.Range("L111").Value = 32999.668
$ws.Range("N111").Value = -41179.668

# Row 116
$ws.Range("H116").Value = 36713.785
$ws.Range("I116").Value = 1090.72
$ws.Range("J116").Value = 333572.66
$ws.Range("K116").Value = 1090.72
$ws.Range("L116").Value = 333572.66
$ws.Range("M116").Value = 1203.28
$ws.Range("N116").Value = -338160.66

# Row 136
$ws.Range("H136").Value = 2786.5
$ws.Range("I136").Value = 1750
$ws.Range("K136").Value = 5250
$ws.Range("M136").Value = -2700

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 36713.785
$ws.Range("I3").Value = 1090.72
$ws.Range("J3").Value = 333572.66
$ws.Range("K3").Value = 1090.72
$ws.Range("L3").Value = 333572.66
$ws.Range("M3").Value = -976.72
$ws.Range("N3").Value = -333800.66

# Row 20
$ws.Range("H20").Value = 25178.395
$ws.Range("I20").Value = 34409.87
$ws.Range("J20").Value = 1330.4166
$ws.Range("K20").Value = 34409.87
$ws.Range("L20").Value = 1330.4166
$ws.Range("M20").Value = -34162.87
$ws.Range("N20").Value = -1824.4166

# Row 132
$ws.Range("H132").Value = 61812.5
$ws.Range("J132").Value = 61812.5
$ws.Range("L132").Value = 61812.5
$ws.Range("N132").Value = -71932.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 6
$ws.Range("H6").Value = 2225.5
$ws.Range("I6").Value = 2000
$ws.Range("J6").Value = 2300.6667
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2300.6667
$ws.Range("M6").Value = -1887
$ws.Range("N6").Value = -2526.6667

# Row 55
$ws.Range("H55").Value = 14178.833
$ws.Range("I55").Value = 10036.5
$ws.Range("J55").Value = 16250
$ws.Range("K55").Value = 10036.5
$ws.Range("L55").Value = 16250
$ws.Range("M55").Value = -9721.5
$ws.Range("N55").Value = -16880

# Row 62
$ws.Range("H62").Value = 2487.5
$ws.Range("I62").Value = 2275
$ws.Range("K62").Value = 2275
$ws.Range("M62").Value = -1651

# Row 65
$ws.Range("H65").Value = 2487.5
$ws.Range("I65").Value = 2275
$ws.Range("K65").Value = 11375
$ws.Range("M65").Value = -8255

# Row 132
$ws.Range("H132").Value = 1995.8223
$ws.Range("I132").Value = 2047.5555
$ws.Range("J132").Value = 1788.8889
$ws.Range("K132").Value = 6142.666499999999
$ws.Range("L132").Value = 5366.6667
$ws.Range("M132").Value = -3612.666499999999
$ws.Range("N132").Value = -10426.6667

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 7
$ws.Range("H7").Value = 522
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 522
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1566
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -1790

# Row 55
$ws.Range("H55").Value = 14440
$ws.Range("J55").Value = 8078.5713
$ws.Range("L55").Value = 24235.7139
$ws.Range("N55").Value = -24589.7139

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 68
$ws.Range("H68").Value = 4053.9
$ws.Range("J68").Value = 4642.375
$ws.Range("L68").Value = 4642.375
$ws.Range("N68").Value = -6140.375

# Row 71
$ws.Range("H71").Value = 4053.9
$ws.Range("J71").Value = 4642.375
$ws.Range("L71").Value = 23211.875
$ws.Range("N71").Value = -30699.875

# Row 81
$ws.Range("H81").Value = 42990
$ws.Range("J81").Value = 42990
$ws.Range("L81").Value = 42990
$ws.Range("N81").Value = -44986

# Row 82
$ws.Range("H82").Value = 1949.9
$ws.Range("I82").Value = 1325.8889
$ws.Range("K82").Value = 1325.8889
$ws.Range("M82").Value = -964.8888999999999

# Row 84
$ws.Range("H84").Value = 42990
$ws.Range("J84").Value = 42990
$ws.Range("L84").Value = 128970
$ws.Range("N84").Value = -138954

# Row 85
$ws.Range("H85").Value = 1949.9
$ws.Range("I85").Value = 1325.8889
$ws.Range("K85").Value = 1325.8889
$ws.Range("M85").Value = -77.88889999999992

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 7
$ws.Range("H7").Value = 1833.3334
$ws.Range("I7").Value = 1500
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -1387
$ws.Range("N7").Value = -2226

# Row 9
$ws.Range("H9").Value = 2000
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

# Row 126
$ws.Range("H126").Value = 1223.5769
$ws.Range("I126").Value = 1286.7368
$ws.Range("J126").Value = 1052.1428
$ws.Range("K126").Value = 3860.2104
$ws.Range("L126").Value = 3156.4284
$ws.Range("M126").Value = -1390.2104
$ws.Range("N126").Value = -8096.428400000001
